# Correccion del archivo Excel de prueba
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo'd "DOCENTEMATERIACODGR" value used across rows 2-4:
# "JuanEdufeedbackAI55555M" -> "Juan - EdufeedbackAI - 55555M"
$ws.Range("D2:D4").Value = "Juan - EdufeedbackAI - 55555M"

# Merge row 5 (2023 / 2023.2) into row 4, keeping row 4's comment,
# then remove the now-duplicate row 5.
$ws.Range("A4").Value = 2023
$ws.Range("B5").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Clear row 5 entirely (its data has been folded into row 4); this
# leaves row 5 empty without shifting row 6 upward.
$ws.Range("A5:J5").ClearContents()

# Add a new "CEDULA" column (K) with header style matching the other
# header cells (G1:J1) and a value of 5 for each data row.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K1").Value = "CEDULA"
$ws.Range("K2").Value = 5
$ws.Range("K3").Value = 5
$ws.Range("K4").Value = 5

# Update the active selection to match the saved view state.
$ws.Range("I17").Select()

$wb.Save()
